$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: period headers
$ws.Range("D8").Value = "9 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "3 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "6 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "9 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "3 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "6 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "9 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# Row 9: publish dates
$ws.Range("D9").Value = "1400-11-05 (3)"
$ws.Range("E9").Value = "1401-04-19 (12)"
$ws.Range("F9").Value = "1401-04-30 (3)"
$ws.Range("G9").Value = "1401-09-15 (4)"
$ws.Range("H9").Value = "1401-11-05 (4)"
$ws.Range("I9").Value = "1402-02-30 (9)"
$ws.Range("J9").Value = "1401-04-30"
$ws.Range("K9").Value = "1401-09-15 (2)"
$ws.Range("L9").Value = "1401-11-05 (2)"
$ws.Range("M9").Value = "1402-02-30"

# Data rows 11-26 (shift left with new M value; row 26 col I also updated)
$ws.Range("D11").Value = 150977
$ws.Range("E11").Value = 244181
$ws.Range("F11").Value = 78509
$ws.Range("G11").Value = 166925
$ws.Range("H11").Value = 268842
$ws.Range("I11").Value = 356254
$ws.Range("J11").Value = 71553
$ws.Range("K11").Value = 150552
$ws.Range("L11").Value = 187805
$ws.Range("M11").Value = 234907
$ws.Range("D12").Value = -67289
$ws.Range("E12").Value = -105352
$ws.Range("F12").Value = -24610
$ws.Range("G12").Value = -63453
$ws.Range("H12").Value = -124418
$ws.Range("I12").Value = -222848
$ws.Range("J12").Value = -48744
$ws.Range("K12").Value = -118307
$ws.Range("L12").Value = -142876
$ws.Range("M12").Value = -171982
$ws.Range("D13").Value = 83688
$ws.Range("E13").Value = 138829
$ws.Range("F13").Value = 53899
$ws.Range("G13").Value = 103471
$ws.Range("H13").Value = 144424
$ws.Range("I13").Value = 133406
$ws.Range("J13").Value = 22809
$ws.Range("K13").Value = 32245
$ws.Range("L13").Value = 44929
$ws.Range("M13").Value = 62924
$ws.Range("D14").Value = -26754
$ws.Range("E14").Value = -35363
$ws.Range("F14").Value = -6563
$ws.Range("G14").Value = -12644
$ws.Range("H14").Value = -17466
$ws.Range("I14").Value = -22677
$ws.Range("J14").Value = -5631
$ws.Range("K14").Value = -13563
$ws.Range("L14").Value = -18895
$ws.Range("M14").Value = -26720
$ws.Range("D16").Value = 17976
$ws.Range("E16").Value = 11559
$ws.Range("F16").Value = 2160
$ws.Range("G16").Value = 1586
$ws.Range("H16").Value = 3180
$ws.Range("I16").Value = 2374
$ws.Range("J16").Value = 2551
$ws.Range("K16").Value = 2854
$ws.Range("L16").Value = 5009
$ws.Range("M16").Value = 369
$ws.Range("D17").Value = 74911
$ws.Range("E17").Value = 115025
$ws.Range("F17").Value = 49496
$ws.Range("G17").Value = 92414
$ws.Range("H17").Value = 130138
$ws.Range("I17").Value = 113103
$ws.Range("J17").Value = 19728
$ws.Range("K17").Value = 21536
$ws.Range("L17").Value = 31044
$ws.Range("M17").Value = 36574
$ws.Range("D18").Value = -8713
$ws.Range("E18").Value = -11862
$ws.Range("F18").Value = -2656
$ws.Range("G18").Value = -4894
$ws.Range("H18").Value = -8192
$ws.Range("I18").Value = -12260
$ws.Range("J18").Value = -4460
$ws.Range("K18").Value = -10643
$ws.Range("L18").Value = -16776
$ws.Range("M18").Value = -22105
$ws.Range("D19").Value = 234567
$ws.Range("E19").Value = 167404
$ws.Range("F19").Value = 1013
$ws.Range("G19").Value = -12632
$ws.Range("H19").Value = -9678
$ws.Range("I19").Value = -6319
$ws.Range("J19").Value = 5542
$ws.Range("K19").Value = 4599
$ws.Range("L19").Value = 22048
$ws.Range("M19").Value = 77647
$ws.Range("D20").Value = 300765
$ws.Range("E20").Value = 270566
$ws.Range("F20").Value = 47853
$ws.Range("G20").Value = 74888
$ws.Range("H20").Value = 112267
$ws.Range("I20").Value = 94524
$ws.Range("J20").Value = 20810
$ws.Range("K20").Value = 15492
$ws.Range("L20").Value = 36316
$ws.Range("M20").Value = 92116
$ws.Range("D21").Value = -3621
$ws.Range("E21").Value = "-"
$ws.Range("F21").Value = "-"
$ws.Range("G21").Value = -16265
$ws.Range("H21").Value = -18771
$ws.Range("I21").Value = -16267
$ws.Range("J21").Value = -2382
$ws.Range("K21").Value = -1568
$ws.Range("L21").Value = -2771
$ws.Range("M21").Value = -2119
$ws.Range("D22").Value = 297144
$ws.Range("E22").Value = 270566
$ws.Range("F22").Value = 47853
$ws.Range("G22").Value = 58622
$ws.Range("H22").Value = 93496
$ws.Range("I22").Value = 78257
$ws.Range("J22").Value = 18428
$ws.Range("K22").Value = 13924
$ws.Range("L22").Value = 33544
$ws.Range("M22").Value = 89997
$ws.Range("D24").Value = 297144
$ws.Range("E24").Value = 270566
$ws.Range("F24").Value = 47853
$ws.Range("G24").Value = 58622
$ws.Range("H24").Value = 93496
$ws.Range("I24").Value = 78257
$ws.Range("J24").Value = 18428
$ws.Range("K24").Value = 13924
$ws.Range("L24").Value = 33544
$ws.Range("M24").Value = 89997
$ws.Range("D26").Value = 4319
$ws.Range("E26").Value = 4202
$ws.Range("F26").Value = 4072
$ws.Range("G26").Value = 3834
$ws.Range("H26").Value = 3652
$ws.Range("I26").Value = 3600
$ws.Range("J26").Value = 105327
$ws.Range("K26").Value = 102502
$ws.Range("L26").Value = 97329
$ws.Range("M26").Value = 87704
